$d = $word.ActiveDocument

# --- 1. First paragraph: append two trailing spaces to the existing text ---
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Microsoft word document.  ", 2)

# --- 2. Append the red "(This is a change – Version for main branch)" text,
#        split across three runs, right before the paragraph mark ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.End = $r.End - 1
$r.Collapse(0)

$r.InsertAfter("(This is a change – Ve")
$r.Font.Color = 255
$r.Collapse(0)

$r.InsertAfter("rsion for main branch")
$r.Font.Color = 255
$r.Collapse(0)

$r.InsertAfter(")")
$r.Font.Color = 255
$r.Collapse(0)

# --- 3. Append a new, empty, shaded paragraph at the very end of the document ---
$d.Content.Find.Execute("we are free at last.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "we are free at last.^p", 2)

$newPara = $d.Paragraphs.Last
$newPara.Range.Style = "Normal"
$newPara.Range.Shading.Texture = 0
$newPara.Range.Shading.ForegroundPatternColor = -16777216
$newPara.Range.Shading.BackgroundPatternColor = 16382457
